$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename application labels in column A (rows 2-15)
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -eq "App1") {
        $cell.Value = "Demo"
    } elseif ($val -eq "App2") {
        $cell.Value = "BeverageStarterFlow"
    } elseif ($val -eq "App3") {
        $cell.Value = "Demo-Brio"
    } elseif ($val -eq "App4") {
        $cell.Value = "Demo-Brio"
    }
}

# Select A15 to match new cursor position
$ws.Range("A15").Select()

$wb.Save()
